# feat: update GenomeEntity Resource
#
# 1) Bump the "Date" property on the Metadata sheet.
# 2) Add a new Concept row ("C0439673" / "Unknown") to the "Include from
#    Unified Medical " sheet's concept table, right above the trailing
#    blank separator row / "System URI" row.

$wb = $excel.ActiveWorkbook

# --- 1) Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2023-04-06T18:12:15+02:00"

# --- 2) Insert a new concept row on the "Include from Unified Medical " sheet ---
$ws = $wb.Worksheets.Item("Include from Unified Medical ")

# Row 8 currently holds the blank separator row and row 9 holds
# "System URI" / the URL. Insert a new row above row 8 so the table grows
# by one row, pushing those two rows down to 9 and 10.
$ws.Rows.Item(8).Insert()

# Insert() hands the new row a blank, un-formatted style; reapply the
# established data-row formatting (border/alignment) from row 7 so the
# new row matches the rest of the table.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(8, 1).Value = "C0439673"
$ws.Cells.Item(8, 2).Value = "Unknown"
